$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 16.97111129760742
$ws.Range("D2").Value = 174

$ws.Range("C3").Value = 16.42203330993652
$ws.Range("D3").Value = 175

$ws.Range("C4").Value = 17.31395721435547
$ws.Range("D4").Value = 176

$ws.Range("C5").Value = 17.78888702392578
$ws.Range("D5").Value = 123

$ws.Range("C6").Value = 16.91699028015137
$ws.Range("D6").Value = 123
